$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "name" header from A1 to I3, clearing A1 in the process.
$ws.Range("I3").Value = $ws.Range("A1").Value2
$ws.Range("A1").ClearContents()
